$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. EnvDetails sheet: update Country (B4) and ProjectName/version value (B6)
# ---------------------------------------------------------------------------
$wsEnv = $wb.Worksheets.Item("EnvDetails")
$wsEnv.Range("B4").Value = "China"
$wsEnv.Range("B6").Value = 1307

# ---------------------------------------------------------------------------
# 2. JIRA_Details sheet: append the new automation scenario rows (17-35)
# ---------------------------------------------------------------------------
$wsJira = $wb.Worksheets.Item("JIRA_Details")

$wsJira.Range("A17").Value = "CreateSubJob"
$wsJira.Range("A18").Value = "PostVendorJournal"
$wsJira.Range("A19").Value = "CreateExpenses"
$wsJira.Range("A20").Value = "Approve_Expenses_Opco"
$wsJira.Range("A21").Value = "Reject_Expenses"
$wsJira.Range("A22").Value = "Post_a_Customer_Payment"
$wsJira.Range("A23").Value = "Customer_Payment_for_Single_Invoice"
$wsJira.Range("A24").Value = "Writing_Off_Bad_Debts"
$wsJira.Range("A25").Value = "Create_Fixed_Asset"
$wsJira.Range("A26").Value = "PostingAssetEntires"
$wsJira.Range("A27").Value = "FixedAssetDisposal"
$wsJira.Range("A28").Value = "FixedAssetReval"
$wsJira.Range("A29").Value = "FixedAssetDepreciation"
$wsJira.Range("A30").Value = "EmployeeCreation"
$wsJira.Range("A31").Value = "EmployeeUserCreation"
$wsJira.Range("A32").Value = "ChangeEmployee"
$wsJira.Range("A33").Value = "CreateUser"
$wsJira.Range("A34").Value = "ChangeUser"
$wsJira.Range("A35").Value = "BlockUser"

$wsJira.Range("A17:A35").NumberFormat = "@"

$wsJira.Range("C17").Value = "TSTAUTO7"
$wsJira.Range("C17").Font.Color = 0

# ---------------------------------------------------------------------------
# 3. View state: EnvDetails is no longer the active tab, JIRA_Details is.
# ---------------------------------------------------------------------------
$wsEnv.Range("B6").Select()

$wsJira.Activate()
$wsJira.Range("B21").Select()
